$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append new row 11 with the latest log entry ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A11").Value = "Status van mijn bestelling"
$logs.Range("B11").Value = "mailmind.test@zohomail.eu"
$logs.Range("C11").Value = "Wanneer wordt mijn bestelling bezorgd?"
$logs.Range("D11").Value = "Bestelling / Levering"
$logs.Range("F11").Value = "2025-06-23 18:24:21"
$logs.Range("G11").Value = "Nee"

# Extend the conditional-formatting ranges so the new row is covered too
# (Modifying one rule's AppliesTo range moves the whole sqref group it belongs to.)
$logs.Range("D2:D10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D11"))
$logs.Range("G2:G10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G11"))

# --- Sheet "Dashboard": re-sort the category/count table now that
#     "Bestelling / Levering" overtook "Offerte / Prijsaanvraag" in count ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Bestelling / Levering"
$dash.Range("A4").Value = "Offerte / Prijsaanvraag"
$dash.Range("B4").Value = 2
